# Reorder the GCI/SD metric rows so that each metric's mix1 rows are
# immediately followed by its mix2 rows (interleaved by metric), instead
# of all mix1 rows followed by all mix2 rows.
#
# Mapping: new row -> old (source) row, for rows 2..37 (row 1 is the header
# and rows stay fixed; rows 2,3,36,37 happen to stay in place as well).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newToOld = @(2,3,20,21,4,5,22,23,6,7,24,25,8,9,26,27,10,11,28,29,12,13,30,31,14,15,32,33,16,17,34,35,18,19,36,37)

$firstRow = 2
$lastRow = 37
$cols = @(1,2,3,4,5,6)   # A..F

# Snapshot all current values for rows 2..37, columns A..F, before we
# overwrite anything (source and destination rows overlap).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($c in $cols) {
        $snapshot[[string]$r + "_" + [string]$c] = $ws.Cells.Item($r, $c).Value2
    }
}

# Write back in the new order.
for ($i = 0; $i -lt $newToOld.Count; $i++) {
    $destRow = $firstRow + $i
    $srcRow = $newToOld[$i]
    foreach ($c in $cols) {
        $key = [string]$srcRow + "_" + [string]$c
        $ws.Cells.Item($destRow, $c).Value2 = $snapshot[$key]
    }
}
